$wb = $excel.ActiveWorkbook

# Rename the second sheet tab
$wsInclude = $wb.Worksheets.Item("Include from Event Types")
$wsInclude.Name = "Include #0"

# Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$wsMeta.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"

# Insert a new row after row 10 (Contact) for the new "Jurisdiction" property.
$wsMeta.Rows.Item(11).Insert()

# Populate the new row 11 with Jurisdiction / empty value.
# A leading apostrophe forces an explicit (empty) text value instead of
# clearing the cell entirely, matching a real empty-string property cell.
$wsMeta.Cells.Item(11, 1).Value = "Jurisdiction"
$wsMeta.Cells.Item(11, 2).Value = "'"

# Re-apply the formatting of the row above it (row 10) so the new row
# matches the style used by the rest of the data rows (the quote-prefix
# write above otherwise leaves its own one-off style behind).
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
